$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1282.6957  # was 1248.44
$ws.Range("J112").Value = 1395.1  # was 1345.9546
$ws.Range("L112").Value = 4185.299999999999  # was 4037.8638
$ws.Range("N112").Value = -6401.299999999999  # was -6253.8638
$ws.Range("H137").Value = 1950.56  # was 1748.4138
$ws.Range("I137").Value = 1464.8572  # was 1245.05
$ws.Range("J137").Value = 2568.7273  # was 2867
$ws.Range("K137").Value = 4394.571599999999  # was 3735.15
$ws.Range("L137").Value = 7706.1819  # was 8601
$ws.Range("M137").Value = -1844.571599999999  # was -1185.15
$ws.Range("N137").Value = -12806.1819  # was -13701
$ws.Range("H138").Value = 2465.0186  # was 2566.818
$ws.Range("I138").Value = 685.5484  # was 698.5333000000001
$ws.Range("J138").Value = 4863.4346  # was 4808.76
$ws.Range("K138").Value = 2056.6452  # was 2095.5999
$ws.Range("L138").Value = 14590.3038  # was 14426.28
$ws.Range("M138").Value = 3083.3548  # was 3044.4001
$ws.Range("N138").Value = -24870.3038  # was -24706.28

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9482.011  # was 13980.85
$ws.Range("I32").Value = 2447.2686  # was 4015.054
$ws.Range("J32").Value = 29974.521  # was 30012.783
$ws.Range("K32").Value = 2447.2686  # was 4015.054
$ws.Range("L32").Value = 29974.521  # was 30012.783
$ws.Range("M32").Value = -2160.2686  # was -3728.054
$ws.Range("N32").Value = -30548.521  # was -30586.783
$ws.Range("H61").Value = 1231.1538  # was 1358.8334
$ws.Range("I61").Value = 1231.1538  # was 1300.5454
$ws.Range("J61").Value = 0  # was 2000
$ws.Range("K61").Value = 1231.1538  # was 1300.5454
$ws.Range("L61").Value = 0  # was 2000
$ws.Range("M61").Value = -1019.1538  # was -1088.5454
$ws.Range("N61").ClearContents()  # was -2424, now removed
$ws.Range("H88").Value = 125001800  # was 66722092
$ws.Range("I88").Value = 1724.75  # was 1679.8
$ws.Range("J88").Value = 250001870  # was 100082296
$ws.Range("K88").Value = 1724.75  # was 1679.8
$ws.Range("L88").Value = 250001870  # was 100082296
$ws.Range("M88").Value = -1318.75  # was -1273.8
$ws.Range("N88").Value = -250002682  # was -100083108
$ws.Range("H91").Value = 125001800  # was 66722092
$ws.Range("I91").Value = 1724.75  # was 1679.8
$ws.Range("J91").Value = 250001870  # was 100082296
$ws.Range("K91").Value = 1724.75  # was 1679.8
$ws.Range("L91").Value = 250001870  # was 100082296
$ws.Range("M91").Value = -320.75  # was -275.8
$ws.Range("N91").Value = -250004678  # was -100085104
$ws.Range("H122").Value = 2749.6  # was 5257
$ws.Range("I122").Value = 995.5  # was 5500
$ws.Range("J122").Value = 3919  # was 5014
$ws.Range("K122").Value = 2986.5  # was 16500
$ws.Range("L122").Value = 11757  # was 15042
$ws.Range("M122").Value = -536.5  # was -14050
$ws.Range("N122").Value = -16657  # was -19942
$ws.Range("H136").Value = 1231.1538  # was 1358.8334
$ws.Range("I136").Value = 1231.1538  # was 1300.5454
$ws.Range("J136").Value = 0  # was 2000
$ws.Range("K136").Value = 3693.4614  # was 3901.6362
$ws.Range("L136").Value = 0  # was 6000
$ws.Range("M136").Value = -1143.4614  # was -1351.6362
$ws.Range("N136").ClearContents()  # was -11100, now removed

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8697694  # was 8335303
$ws.Range("I86").Value = 12501910  # was 11766562
$ws.Range("J86").Value = 2342.8572  # was 2243.8572
$ws.Range("K86").Value = 12501910  # was 11766562
$ws.Range("L86").Value = 2342.8572  # was 2243.8572
$ws.Range("M86").Value = -12500787  # was -11765439
$ws.Range("N86").Value = -4588.8572  # was -4489.8572
$ws.Range("H89").Value = 8697694  # was 8335303
$ws.Range("I89").Value = 12501910  # was 11766562
$ws.Range("J89").Value = 2342.8572  # was 2243.8572
$ws.Range("K89").Value = 62509550  # was 58832810
$ws.Range("L89").Value = 11714.286  # was 11219.286
$ws.Range("M89").Value = -62503934  # was -58827194
$ws.Range("N89").Value = -22946.286  # was -22451.286
$ws.Range("H134").Value = 2907.5557  # was 2884.111
$ws.Range("I134").Value = 2404.6  # was 2362.4
$ws.Range("K134").Value = 7213.799999999999  # was 7087.200000000001
$ws.Range("M134").Value = -4678.799999999999  # was -4552.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6485255.5  # was 6701441
$ws.Range("I31").Value = 8740223  # was 9572544
$ws.Range("J31").Value = 2225  # was 2200
$ws.Range("K31").Value = 8740223  # was 9572544
$ws.Range("L31").Value = 2225  # was 2200
$ws.Range("M31").Value = -8739928  # was -9572249
$ws.Range("N31").Value = -2815  # was -2790
$ws.Range("H34").Value = 6485255.5  # was 6701441
$ws.Range("I34").Value = 8740223  # was 9572544
$ws.Range("J34").Value = 2225  # was 2200
$ws.Range("K34").Value = 8740223  # was 9572544
$ws.Range("L34").Value = 2225  # was 2200
$ws.Range("M34").Value = -8740021  # was -9572342
$ws.Range("N34").Value = -2629  # was -2604
$ws.Range("H58").Value = 1918.4736  # was 2161.3125
$ws.Range("I58").Value = 1175  # was 1526.6666
$ws.Range("J58").Value = 2261.6155  # was 2307.7693
$ws.Range("K58").Value = 1175  # was 1526.6666
$ws.Range("L58").Value = 2261.6155  # was 2307.7693
$ws.Range("M58").Value = -972  # was -1323.6666
$ws.Range("N58").Value = -2667.6155  # was -2713.7693
$ws.Range("H132").Value = 1297.3793  # was 1914
$ws.Range("I132").Value = 963.5833  # was 1309.3334
$ws.Range("J132").Value = 2899.6  # was 4332.6665
$ws.Range("K132").Value = 2890.7499  # was 3928.0002
$ws.Range("L132").Value = 8698.799999999999  # was 12997.9995
$ws.Range("M132").Value = -360.7498999999998  # was -1398.0002
$ws.Range("N132").Value = -13758.8  # was -18057.9995
$ws.Range("H136").Value = 1918.4736  # was 2161.3125
$ws.Range("I136").Value = 1175  # was 1526.6666
$ws.Range("J136").Value = 2261.6155  # was 2307.7693
$ws.Range("K136").Value = 3525  # was 4579.9998
$ws.Range("L136").Value = 6784.8465  # was 6923.3079
$ws.Range("M136").Value = -975  # was -2029.9998
$ws.Range("N136").Value = -11884.8465  # was -12023.3079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 914.5238000000001  # was 903.2368
$ws.Range("I131").Value = 423.6  # was 340
$ws.Range("J131").Value = 980.86487  # was 988.57574
$ws.Range("K131").Value = 1270.8  # was 1020
$ws.Range("L131").Value = 2942.59461  # was 2965.72722
$ws.Range("M131").Value = 3769.2  # was 4020
$ws.Range("N131").Value = -13022.59461  # was -13045.72722
$ws.Range("H137").Value = 5581.5186  # was 6382.909
$ws.Range("I137").Value = 966.6667  # was 671.7273
$ws.Range("J137").Value = 6158.375  # was 12094.091
$ws.Range("K137").Value = 2900.0001  # was 2015.1819
$ws.Range("L137").Value = 18475.125  # was 36282.273
$ws.Range("M137").Value = 2199.9999  # was 3084.8181
$ws.Range("N137").Value = -28675.125  # was -46482.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 33338332  # was 25810994
$ws.Range("I70").Value = 66670924  # was 57146708
$ws.Range("J70").Value = 5740  # was 5110.5293
$ws.Range("K70").Value = 66670924  # was 57146708
$ws.Range("L70").Value = 5740  # was 5110.5293
$ws.Range("M70").Value = -66670654  # was -57146438
$ws.Range("N70").Value = -6280  # was -5650.5293
$ws.Range("H73").Value = 33338332  # was 25810994
$ws.Range("I73").Value = 66670924  # was 57146708
$ws.Range("J73").Value = 5740  # was 5110.5293
$ws.Range("K73").Value = 66670924  # was 57146708
$ws.Range("L73").Value = 5740  # was 5110.5293
$ws.Range("M73").Value = -66669988  # was -57145772
$ws.Range("N73").Value = -7612  # was -6982.5293
$ws.Range("H132").Value = 2221.8635  # was 2200.6538
$ws.Range("I132").Value = 1492.2  # was 1484.3889
$ws.Range("J132").Value = 3785.4285  # was 3812.25
$ws.Range("K132").Value = 4476.6  # was 4453.1667
$ws.Range("L132").Value = 11356.2855  # was 11436.75
$ws.Range("M132").Value = -1946.6  # was -1923.1667
$ws.Range("N132").Value = -16416.2855  # was -16496.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 384.32144  # was 629.2
$ws.Range("I16").Value = 405.03845  # was 629.2
$ws.Range("J16").Value = 115  # was 0
$ws.Range("K16").Value = 405.03845  # was 629.2
$ws.Range("L16").Value = 115  # was 0
$ws.Range("M16").Value = -235.03845  # was -459.2
$ws.Range("N16").Value = -455  # new cell
$ws.Range("H46").Value = 825.88  # was 779.6070999999999
$ws.Range("I46").Value = 629  # was 607.7143
$ws.Range("J46").Value = 888.0526  # was 836.9048
$ws.Range("K46").Value = 629  # was 607.7143
$ws.Range("L46").Value = 888.0526  # was 836.9048
$ws.Range("M46").Value = -441  # was -419.7143
$ws.Range("N46").Value = -1264.0526  # was -1212.9048
